$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 309, shifting existing rows 309:382 down to 310:383
$ws.Rows.Item(309).Insert()

# Populate the newly inserted row 309 with the new record's data
$ws.Cells.Item(309, 1).Value = 4
$ws.Cells.Item(309, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(309, 3).Value = "Los Lagos"
$ws.Cells.Item(309, 4).Value = 44932
$ws.Cells.Item(309, 5).Value = 10
$ws.Cells.Item(309, 6).Value = 100112040
$ws.Cells.Item(309, 7).Value = "Cilantro"
$ws.Cells.Item(309, 8).Value = "Sin especificar"
$ws.Cells.Item(309, 9).Value = "Primera"
$ws.Cells.Item(309, 10).Value = 120
$ws.Cells.Item(309, 11).Value = 10000
$ws.Cells.Item(309, 12).Value = 10000
$ws.Cells.Item(309, 13).Value = 10000
$ws.Cells.Item(309, 14).Value = "$/docena de atados (2 kilos)"
$ws.Cells.Item(309, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(309, 16).Value = 5000
$ws.Cells.Item(309, 17).Value = 2
$ws.Cells.Item(309, 18).Value = "Hortaliza"
